$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Escopo")

# Row 7 - item "O aluno criou as migrações do banco de dados?" - porcentagem 0.25 -> 0.5
$ws.Range("F7").Value = 0.5

# Row 8 - item "O aluno criou as operações de CRUD?" - status "Em Andamento" -> "Concluido"
$ws.Range("E8").Value = "Concluido"

# Row 23 - item "O aluno desenvolveu um microserviço CRUD?" - status "Em Andamento" -> "Concluido", porcentagem 0.5 -> 1
$ws.Range("E23").Value = "Concluido"
$ws.Range("F23").Value = 1

# Row 25 - item "O sistema implementado funcionou corretamente?" - status blank -> "Em Andamento", porcentagem 0 -> 0.1
$ws.Range("E25").Value = "Em Andamento"
$ws.Range("F25").Value = 0.1

# Update the active selection to D11 as recorded in the sheet view
$ws.Range("D11").Select() | Out-Null
